$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update USERID cell (G2): 32118 -> 31160
$ws.Range("G2").Value = 31160

# Update PREPARATION cell (F2) text to reflect the new username
$ws.Range("F2").Value = "Username : 31160;" + [char]10 + "Password : bni1234;" + [char]10 + "Tgl. Market : 23/01/2023;" + [char]10 + "File Excel : 23012023HargaPasarFixedIncome.xlsx"

# Update sheet view: scroll/selection moved from column G to column F/E
$win = $excel.ActiveWindow
$ws.Range("F3").Select() | Out-Null
$win.ScrollColumn = 5
$win.ScrollRow = 1
